$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New breakdown rows (appended after the existing 3 rows of data)
$rows = @(
    @("2025-02-07", "V3", "Blade Tension"),
    @("2025-02-07", "Rivers", "Oil Pressure"),
    @("2025-02-07", "K4", "Oil Spring Problem"),
    @("2025-02-07", "J1", "Blade Tension"),
    @("2025-02-07", "J1", "Blade Tension"),
    @("2025-02-07", "J1", "Blade Tension"),
    @("2025-02-07", "ITM2", "Blade Guide Problem"),
    @("2025-02-07", "K4", "Blade Tension")
)

$r = 4
foreach ($row in $rows) {
    # Leading apostrophe forces the date-looking text to stay text (quotePrefix),
    # matching the existing rows which store the date as plain text too.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Column widths (values chosen so the engine's internal rounding lands as
# close as possible to the target stored widths of 29 / 28.90625 / 30.453125)
$ws.Columns.Item(1).ColumnWidth = 28.16666667
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(3).ColumnWidth = 29.66666667

# Move the active selection to B1
$ws.Range("B1").Select() | Out-Null
